# longlist update, shortlist mobile banking finished
#
# 1. Add new "Anmerkung" (remark) texts in column N of the "Mobile Banking"
#    sheet for each provider (longlist update).
# 2. Widen column N to fit the new remarks.
# 3. Make "Mobile Banking" (sheet 1) the active sheet / selected cell,
#    instead of "Contactless Payment" (shortlist finished there).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mobile Banking")

# New remarks, entered in the same order the original author must have
# typed them (this drives the order new shared-string entries are created
# in, and therefore which index each one receives).
$remarks = @(
  @{Cell = "N3";  Text = "gegründet 1968, weltweit über 14000 Kunden"},
  @{Cell = "N5";  Text = "gegründet 1984, unterstützt 1000 Institutionen weltweit"},
  @{Cell = "N4";  Text = "gegründet 2003, Dienstleister von 350 Institutionen weltweit"},
  @{Cell = "N6";  Text = "gegründet 1981, Hauptsitz in Indien"},
  @{Cell = "N7";  Text = "gegründet 1995"},
  @{Cell = "N8";  Text = "gegründet 1976, Partner von über 1300 Banken"},
  @{Cell = "N9";  Text = "gegründet 2007,  über 20Mil. Nutzer von über 600 Apps"},
  @{Cell = "N10"; Text = "Tochter der SAP-Gruppe"},
  @{Cell = "N11"; Text = "gegründet 1975, Dienstleister von 1650 Banken und Institutionen"},
  @{Cell = "N13"; Text = "gegründet 1963, weltweit über 1800 Kunden"},
  @{Cell = "N14"; Text = "gegründet 2005"}
)

foreach ($remark in $remarks) {
    $cell = $ws.Range($remark.Cell)
    $cell.Value = $remark.Text
    # match the left-aligned "General" style used throughout column A/N
    $cell.HorizontalAlignment = -4131
}

# Rows without a remark still received the same formatting (empty cell
# with the style applied) in column N.
foreach ($ref in @("N2", "N12", "N15")) {
    $ws.Range($ref).HorizontalAlignment = -4131
}

# Widen column N so the longer remarks are readable.
$ws.Columns.Item(14).ColumnWidth = 57.14

# "Mobile Banking" becomes the active sheet, scrolled/selected near the
# new column, replacing "Contactless Payment" as the active tab.
$ws.Activate() | Out-Null
$ws.Range("G1").Select() | Out-Null
$ws.Range("N14").Select() | Out-Null
